$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.727.50'
$ws.Range('E2').Value = '  -2.38%  '
$ws.Range('D3').Value = '3.497.15'
$ws.Range('E3').Value = '  -3.71%  '
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = "'605.37"
$ws.Range('E5').Value = '  -3.33%  '
$ws.Range('D6').Value = "'150.30"
$ws.Range('E6').Value = '  -5.36%  '
$ws.Range('D7').Value = '3.495.51'
$ws.Range('E7').Value = '  -3.77%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('E10').Value = '  +0.01%  '
$ws.Range('D11').Value = "'7.53"
$ws.Range('E11').Value = '  +4.72%  '
$ws.Range('E12').Value = '  -1.75%  '
$ws.Range('E13').Value = '  -4.23%  '
$ws.Range('D14').Value = "'31.97"
$ws.Range('E14').Value = '  -3.96%  '
$ws.Range('D15').Value = '4.089.58'
$ws.Range('E15').Value = '  -3.56%  '
$ws.Range('D16').Value = '67.746.60'
$ws.Range('E16').Value = '  -1.79%  '
$ws.Range('D17').Value = '3.493.77'
$ws.Range('E17').Value = '  -3.21%  '
$ws.Range('D18').Value = "'0.117"
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').Value = "'6.50"
$ws.Range('E19').Value = '  -0.96%  '
$ws.Range('D20').Value = "'15.34"
$ws.Range('E20').Value = '  -2.48%  '
$ws.Range('D21').Value = "'9.94"
$ws.Range('E21').Value = '  -2.08%  '
$ws.Range('D22').Value = "'445.02"
$ws.Range('E22').Value = '  -4.31%  '
$ws.Range('D23').Value = "'0.625"
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('D24').Value = "'78.92"
$ws.Range('E24').Value = '  +0.50%  '
$ws.Range('D25').Value = '3.640.55'
$ws.Range('E25').Value = '  -3.42%  '
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('E27').Value = '  -8.90%  '
$ws.Range('D28').Value = "'8.67"
$ws.Range('E28').Value = '  -4.98%  '
$ws.Range('D29').Value = "'9.93"
$ws.Range('E29').Value = '  -5.55%  '
$ws.Range('D30').Value = "'1.66"
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('E31').Value = '  -4.23%  '
$ws.Range('E32').Value = '  -2.92%  '
$ws.Range('E33').Value = '  +0.34%  '
$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').Value = "'25.57"
$ws.Range('E34').Value = '  -3.32%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = "'6.17"
$ws.Range('E35').Value = '  -5.81%  '
$ws.Range('E36').Value = '  -4.99%  '
$ws.Range('D37').Value = '3.491.34'
$ws.Range('E37').Value = '  -3.40%  '
$ws.Range('E38').Value = '  -4.10%  '
$ws.Range('E39').Value = '  -0.02%  '
$ws.Range('D40').Value = "'2.31"
$ws.Range('E40').Value = '  -2.04%  '
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = "'176.46"
$ws.Range('E42').Value = '  +0.21%  '
$ws.Range('E43').Value = '  -2.93%  '
$ws.Range('E44').Value = '  -3.72%  '
$ws.Range('D45').Value = "'0.898"
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').Value = "'29.94"
$ws.Range('E46').Value = '  -1.14%  '
$ws.Range('D47').Value = "'46.90"
$ws.Range('E47').Value = '  +1.75%  '
$ws.Range('D48').Value = "'1.29"
$ws.Range('E48').Value = '  -5.17%  '
$ws.Range('D49').Value = "'2.51"
$ws.Range('E49').Value = '  -9.69%  '
$ws.Range('D50').Value = "'7.61"
$ws.Range('E50').Value = '  -1.96%  '
$ws.Range('D51').Value = "'0.995"
$ws.Range('E51').Value = '  -3.26%  '
